$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (previously row "44875" entry) -> now takes the old row 3 values for D,M,N,O,P,Q,S,T (R unchanged)
$ws.Range("D2").Value = 44855
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = "$/bandeja 5 kilos"
$ws.Range("S2").Value = 3000
$ws.Range("T2").Value = 5

# Row 3 -> now takes the old row 5 values for D,M,Q,R,S,T (N,O,P unchanged)
$ws.Range("D3").Value = 45222
$ws.Range("M3").Value = 20
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 1500
$ws.Range("T3").Value = 10

# Row 4 -> now takes the old row 2 values for D,M (N,O,P,Q,R,S,T unchanged)
$ws.Range("D4").Value = 44875
$ws.Range("M4").Value = 50

# Row 5 -> now takes the old row 4 values for D,M,N,O,P,R,S (Q,T unchanged)
$ws.Range("D5").Value = 44874
$ws.Range("M5").Value = 67
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 16000
$ws.Range("R5").Value = "Provincia de Los Andes"
$ws.Range("S5").Value = 1600
